$d = $word.ActiveDocument

# Locate the target paragraph ("Créer la méthode de classe
# trouver_personnes_majeurs() qui retourne ...") by its distinctive
# (pre-edit) text, then expand the found range so it covers the whole
# paragraph (including its end-of-paragraph mark).
$rng = $d.Content
$ok = $rng.Find.Execute("trouver_personnes_majeurs() qui retourne la liste des personnes qui ont plus de 18 ans.")
if (-not $ok) {
    throw "Target paragraph not found"
}
$rng.Expand(4)  # wdParagraph

$xml = '<w:p w:rsidR="00924AC8" w:rsidRDefault="00924AC8" w:rsidP="00924AC8"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/></w:rPr><w:t xml:space="preserve">Créer la méthode de classe </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="0D0D0D"/></w:rPr><w:t>trouver_personnes_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="0D0D0D"/></w:rPr><w:t>majeurs(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="0D0D0D"/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="0D0D0D"/></w:rPr><w:t xml:space="preserve"> qui retourne la liste des personnes qui ont plus de 18 ans.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:color w:val="auto"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>Utiliser un dictionnaire statique pour sauvegarder les instances de la classe Personne. La clé du dictionnaire c’est le n</w:t></w:r><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>uméro d’assurance sociale</w:t></w:r><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve"> et la valeur c’est l’</w:t></w:r><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>instance.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:color w:val="auto"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="auto"/></w:rPr><w:t>Tests unitaires avec Pytest</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

# Replace the whole paragraph's contents (raw OOXML) with: the same
# paragraph reworded ("trouver_personnes_" / "majeurs(" / ")" as three
# runs, with <w:proofErr> grammar markers around "majeurs("), followed
# by the two new list items ("Utiliser un dictionnaire statique ..."
# and "Tests unitaires avec Pytest") — the _GoBack bookmark now sits at
# the end of the last of these new paragraphs.
$rng.InsertXML($xml)

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
